$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 27 de Julio de 2020 a las 02:04'

# Row 4
$ws.Cells.Item(4, 1).Value = 'Estados Unidos'
$ws.Cells.Item(4, 2).Value = 4370522
$ws.Cells.Item(4, 3).Value = 54813
$ws.Cells.Item(4, 4).Value = 2087667
$ws.Cells.Item(4, 5).Value = 2133029
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 428
$ws.Cells.Item(4, 8).Value = 149826

# Row 5
$ws.Cells.Item(5, 1).Value = 'Brasil'
$ws.Cells.Item(5, 2).Value = 2419901
$ws.Cells.Item(5, 3).Value = 23467
$ws.Cells.Item(5, 4).Value = 1634274
$ws.Cells.Item(5, 5).Value = 698575
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 556
$ws.Cells.Item(5, 8).Value = 87052

# Row 17
$ws.Cells.Item(17, 1).Value = 'Colombia'
$ws.Cells.Item(17, 2).Value = 248976
$ws.Cells.Item(17, 3).Value = 8181
$ws.Cells.Item(17, 4).Value = 125037
$ws.Cells.Item(17, 5).Value = 115414
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 256
$ws.Cells.Item(17, 8).Value = 8525

# Row 18
$ws.Cells.Item(18, 1).Value = 'Italia'
$ws.Cells.Item(18, 2).Value = 246118
$ws.Cells.Item(18, 3).Value = 254
$ws.Cells.Item(18, 4).Value = 198446
$ws.Cells.Item(18, 5).Value = 12565
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(18, 7).Value = 5
$ws.Cells.Item(18, 8).Value = 35107

# Row 24
$ws.Cells.Item(24, 1).Value = 'Canada'
$ws.Cells.Item(24, 2).Value = 113911
$ws.Cells.Item(24, 3).Value = 355
$ws.Cells.Item(24, 4).Value = 99355
$ws.Cells.Item(24, 5).Value = 5666
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(24, 7).Value = 5
$ws.Cells.Item(24, 8).Value = 8890

# Row 56
$ws.Cells.Item(56, 1).Value = 'Ghana'
$ws.Cells.Item(56, 2).Value = 32969
$ws.Cells.Item(56, 3).Value = 1118
$ws.Cells.Item(56, 4).Value = 29494
$ws.Cells.Item(56, 5).Value = 3307
$ws.Cells.Item(56, 6).Value = 0
$ws.Cells.Item(56, 7).Value = 7
$ws.Cells.Item(56, 8).Value = 168

# Row 57
$ws.Cells.Item(57, 1).Value = 'Kirguistan'
$ws.Cells.Item(57, 2).Value = 32813
$ws.Cells.Item(57, 3).Value = 689
$ws.Cells.Item(57, 4).Value = 20388
$ws.Cells.Item(57, 5).Value = 11148
$ws.Cells.Item(57, 6).Value = 0
$ws.Cells.Item(57, 7).Value = 28
$ws.Cells.Item(57, 8).Value = 1277

# Row 71
$ws.Cells.Item(71, 1).Value = 'Venezuela'
$ws.Cells.Item(71, 2).Value = 15463
$ws.Cells.Item(71, 3).Value = 534
$ws.Cells.Item(71, 4).Value = 9746
$ws.Cells.Item(71, 5).Value = 5575
$ws.Cells.Item(71, 6).Value = 0
$ws.Cells.Item(71, 7).Value = 4
$ws.Cells.Item(71, 8).Value = 142

# Row 72
$ws.Cells.Item(72, 1).Value = 'Chequia'
$ws.Cells.Item(72, 2).Value = 15324
$ws.Cells.Item(72, 3).Value = 112
$ws.Cells.Item(72, 4).Value = 11428
$ws.Cells.Item(72, 5).Value = 3525
$ws.Cells.Item(72, 6).Value = 0
$ws.Cells.Item(72, 7).Value = 2
$ws.Cells.Item(72, 8).Value = 371

# Row 73
$ws.Cells.Item(73, 1).Value = 'Costa Rica'
$ws.Cells.Item(73, 2).Value = 15229
$ws.Cells.Item(73, 3).Value = 629
$ws.Cells.Item(73, 4).Value = 3736
$ws.Cells.Item(73, 5).Value = 11389
$ws.Cells.Item(73, 6).Value = 0
$ws.Cells.Item(73, 7).Value = 6
$ws.Cells.Item(73, 8).Value = 104

# Row 86
$ws.Cells.Item(86, 1).Value = 'Noruega'
$ws.Cells.Item(86, 2).Value = 9117
$ws.Cells.Item(86, 3).Value = 6
$ws.Cells.Item(86, 4).Value = 8752
$ws.Cells.Item(86, 5).Value = 110
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(86, 8).Value = 255

# Row 93
$ws.Cells.Item(93, 1).Value = 'Guinea'
$ws.Cells.Item(93, 2).Value = 7008
$ws.Cells.Item(93, 3).Value = 81
$ws.Cells.Item(93, 4).Value = 6152
$ws.Cells.Item(93, 5).Value = 813
$ws.Cells.Item(93, 6).Value = 0
$ws.Cells.Item(93, 7).Value = 1
$ws.Cells.Item(93, 8).Value = 43

# Row 94
$ws.Cells.Item(94, 1).Value = 'Gabon'
$ws.Cells.Item(94, 2).Value = 6984
$ws.Cells.Item(94, 3).Value = 0
$ws.Cells.Item(94, 4).Value = 4463
$ws.Cells.Item(94, 5).Value = 2472
$ws.Cells.Item(94, 6).Value = 0
$ws.Cells.Item(94, 7).Value = 0
$ws.Cells.Item(94, 8).Value = 49

# Row 95
$ws.Cells.Item(95, 1).Value = 'Luxemburgo'
$ws.Cells.Item(95, 2).Value = 6272
$ws.Cells.Item(95, 3).Value = 83
$ws.Cells.Item(95, 4).Value = 4647
$ws.Cells.Item(95, 5).Value = 1513
$ws.Cells.Item(95, 6).Value = 0
$ws.Cells.Item(95, 7).Value = 0
$ws.Cells.Item(95, 8).Value = 112

# Row 96
$ws.Cells.Item(96, 1).Value = 'Mauritania'
$ws.Cells.Item(96, 2).Value = 6171
$ws.Cells.Item(96, 3).Value = 20
$ws.Cells.Item(96, 4).Value = 4430
$ws.Cells.Item(96, 5).Value = 1585
$ws.Cells.Item(96, 6).Value = 0
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(96, 8).Value = 156

# Row 97
$ws.Cells.Item(97, 1).Value = 'Republica de Yibuti'
$ws.Cells.Item(97, 2).Value = 5050
$ws.Cells.Item(97, 3).Value = 11
$ws.Cells.Item(97, 4).Value = 4966
$ws.Cells.Item(97, 5).Value = 26
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 58

# Row 100
$ws.Cells.Item(100, 1).Value = 'Republica de Africa Central'
$ws.Cells.Item(100, 2).Value = 4599
$ws.Cells.Item(100, 3).Value = 1
$ws.Cells.Item(100, 4).Value = 1546
$ws.Cells.Item(100, 5).Value = 2994
$ws.Cells.Item(100, 6).Value = 0
$ws.Cells.Item(100, 7).Value = 0
$ws.Cells.Item(100, 8).Value = 59

# Row 102
$ws.Cells.Item(102, 1).Value = 'Paraguay'
$ws.Cells.Item(102, 2).Value = 4444
$ws.Cells.Item(102, 3).Value = 116
$ws.Cells.Item(102, 4).Value = 2794
$ws.Cells.Item(102, 5).Value = 1609
$ws.Cells.Item(102, 6).Value = 0
$ws.Cells.Item(102, 7).Value = 1
$ws.Cells.Item(102, 8).Value = 41

# Row 103
$ws.Cells.Item(103, 1).Value = 'Hungria'
$ws.Cells.Item(103, 2).Value = 4435
$ws.Cells.Item(103, 3).Value = 11
$ws.Cells.Item(103, 4).Value = 3329
$ws.Cells.Item(103, 5).Value = 510
$ws.Cells.Item(103, 6).Value = 0
$ws.Cells.Item(103, 7).Value = 0
$ws.Cells.Item(103, 8).Value = 596

# Row 106
$ws.Cells.Item(106, 1).Value = 'Malaui'
$ws.Cells.Item(106, 2).Value = 3640
$ws.Cells.Item(106, 3).Value = 83
$ws.Cells.Item(106, 4).Value = 1639
$ws.Cells.Item(106, 5).Value = 1902
$ws.Cells.Item(106, 6).Value = 0
$ws.Cells.Item(106, 7).Value = 5
$ws.Cells.Item(106, 8).Value = 99

# Row 114
$ws.Cells.Item(114, 1).Value = 'Montenegro'
$ws.Cells.Item(114, 2).Value = 2799
$ws.Cells.Item(114, 3).Value = 52
$ws.Cells.Item(114, 4).Value = 739
$ws.Cells.Item(114, 5).Value = 2017
$ws.Cells.Item(114, 6).Value = 0
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 43

# Row 115
$ws.Cells.Item(115, 1).Value = 'Sri Lanka'
$ws.Cells.Item(115, 2).Value = 2782
$ws.Cells.Item(115, 3).Value = 12
$ws.Cells.Item(115, 4).Value = 2106
$ws.Cells.Item(115, 5).Value = 665
$ws.Cells.Item(115, 6).Value = 0
$ws.Cells.Item(115, 7).Value = 0
$ws.Cells.Item(115, 8).Value = 11

# Row 118
$ws.Cells.Item(118, 1).Value = 'Zimbabue'
$ws.Cells.Item(118, 2).Value = 2512
$ws.Cells.Item(118, 3).Value = 78
$ws.Cells.Item(118, 4).Value = 518
$ws.Cells.Item(118, 5).Value = 1960
$ws.Cells.Item(118, 6).Value = 0
$ws.Cells.Item(118, 7).Value = 0
$ws.Cells.Item(118, 8).Value = 34

# Row 119
$ws.Cells.Item(119, 1).Value = 'Mali'
$ws.Cells.Item(119, 2).Value = 2510
$ws.Cells.Item(119, 3).Value = 7
$ws.Cells.Item(119, 4).Value = 1911
$ws.Cells.Item(119, 5).Value = 476
$ws.Cells.Item(119, 6).Value = 0
$ws.Cells.Item(119, 7).Value = 0
$ws.Cells.Item(119, 8).Value = 123

# Row 120
$ws.Cells.Item(120, 1).Value = 'Cuba'
$ws.Cells.Item(120, 2).Value = 2495
$ws.Cells.Item(120, 3).Value = 17
$ws.Cells.Item(120, 4).Value = 2349
$ws.Cells.Item(120, 5).Value = 59
$ws.Cells.Item(120, 6).Value = 0
$ws.Cells.Item(120, 7).Value = 0
$ws.Cells.Item(120, 8).Value = 87

# Row 121
$ws.Cells.Item(121, 1).Value = 'Cabo Verde'
$ws.Cells.Item(121, 2).Value = 2307
$ws.Cells.Item(121, 3).Value = 49
$ws.Cells.Item(121, 4).Value = 1447
$ws.Cells.Item(121, 5).Value = 838
$ws.Cells.Item(121, 6).Value = 0
$ws.Cells.Item(121, 7).Value = 0
$ws.Cells.Item(121, 8).Value = 22

# Row 122
$ws.Cells.Item(122, 1).Value = 'Sudan del Sur'
$ws.Cells.Item(122, 2).Value = 2262
$ws.Cells.Item(122, 3).Value = 4
$ws.Cells.Item(122, 4).Value = 1175
$ws.Cells.Item(122, 5).Value = 1042
$ws.Cells.Item(122, 6).Value = 0
$ws.Cells.Item(122, 7).Value = 0
$ws.Cells.Item(122, 8).Value = 45

# Row 138
$ws.Cells.Item(138, 1).Value = 'Surinam'
$ws.Cells.Item(138, 2).Value = 1439
$ws.Cells.Item(138, 3).Value = 58
$ws.Cells.Item(138, 4).Value = 890
$ws.Cells.Item(138, 5).Value = 526
$ws.Cells.Item(138, 6).Value = 0
$ws.Cells.Item(138, 7).Value = 0
$ws.Cells.Item(138, 8).Value = 23

# Row 140
$ws.Cells.Item(140, 1).Value = 'Uruguay'
$ws.Cells.Item(140, 2).Value = 1192
$ws.Cells.Item(140, 3).Value = 18
$ws.Cells.Item(140, 4).Value = 948
$ws.Cells.Item(140, 5).Value = 210
$ws.Cells.Item(140, 6).Value = 0
$ws.Cells.Item(140, 7).Value = 0
$ws.Cells.Item(140, 8).Value = 34

# Row 143
$ws.Cells.Item(143, 1).Value = 'Niger'
$ws.Cells.Item(143, 2).Value = 1136
$ws.Cells.Item(143, 3).Value = 12
$ws.Cells.Item(143, 4).Value = 1027
$ws.Cells.Item(143, 5).Value = 40
$ws.Cells.Item(143, 6).Value = 0
$ws.Cells.Item(143, 7).Value = 0
$ws.Cells.Item(143, 8).Value = 69

# Row 144
$ws.Cells.Item(144, 1).Value = 'Georgia'
$ws.Cells.Item(144, 2).Value = 1131
$ws.Cells.Item(144, 3).Value = 14
$ws.Cells.Item(144, 4).Value = 920
$ws.Cells.Item(144, 5).Value = 195
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 16

# Row 147
$ws.Cells.Item(147, 1).Value = 'Republica de Chipre'
$ws.Cells.Item(147, 2).Value = 1057
$ws.Cells.Item(147, 3).Value = 4
$ws.Cells.Item(147, 4).Value = 852
$ws.Cells.Item(147, 5).Value = 186
$ws.Cells.Item(147, 6).Value = 0
$ws.Cells.Item(147, 7).Value = 0
$ws.Cells.Item(147, 8).Value = 19

# Row 169
$ws.Cells.Item(169, 1).Value = 'Bahamas'
$ws.Cells.Item(169, 2).Value = 342
$ws.Cells.Item(169, 3).Value = 16
$ws.Cells.Item(169, 4).Value = 91
$ws.Cells.Item(169, 5).Value = 240
$ws.Cells.Item(169, 6).Value = 0
$ws.Cells.Item(169, 7).Value = 0
$ws.Cells.Item(169, 8).Value = 11

# Row 170
$ws.Cells.Item(170, 1).Value = 'Isla de Man'
$ws.Cells.Item(170, 2).Value = 336
$ws.Cells.Item(170, 3).Value = 0
$ws.Cells.Item(170, 4).Value = 312
$ws.Cells.Item(170, 5).Value = 0
$ws.Cells.Item(170, 6).Value = 0
$ws.Cells.Item(170, 7).Value = 0
$ws.Cells.Item(170, 8).Value = 24

